$wb = $excel.ActiveWorkbook

# Restricciones_del_follower sheet
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-32.425 + x + 6.833333333333332y"
$ws.Range("B2").Value = "'18.424999999999994"
$ws.Range("D2").Value = "'0.07"
$ws.Range("E2").Value = "'6.0"
$ws.Range("F2").Value = "'4.1"
$ws.Range("A3").Value = "-1.5100000000000007 + x - 0.7999999999999998y"
$ws.Range("B3").Value = "'-0.4899999999999993"
$ws.Range("D3").Value = "'0.21"
$ws.Range("E3").Value = "'0.6"
$ws.Range("F3").Value = "'7.1"
$ws.Range("A4").Value = "-6.775000000000001 - 2x + 1.1666666666666665y"
$ws.Range("B4").Value = "'-5.775000000000001"
$ws.Range("D4").Value = "'0.4"
$ws.Range("E4").Value = "'6.1"
$ws.Range("F4").Value = "'0.7000000000000001"

# Punto_modificado sheet
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "'4.75"
$ws.Range("B2").Value = "'4.05"

# Vector_bf sheet
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "'1.1230000000000004"

# Vector_BF sheet (index 6; name lookup is case-insensitive and would
# otherwise collide with "Vector_bf", sheet 5)
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = "'2.0999999999999996"
$ws.Range("A3").Value = "'-52.93666666666666"

# Vector_Alpha sheet
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.6000000000000001
